$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (header stays at row 1; rows 2..236 hold the old data)
$ws.Range("A2:C236").EntireRow.Delete()

# Make room for the new (re-sorted / edited) data block: 236 rows
$ws.Range("A2:A237").EntireRow.Insert()

$data = @(
  ,("005103059","Walquiria",25855.71)
  ,("005142611","Guilherme",22441.24)
  ,("005531186","Rafael",9912)
  ,("004212438","Kenia",3750.01)
  ,("004368468","Ahmad",1966.84)
  ,("004332747","Lohran",989.79)
  ,("004392159","Rodrigo",902.31)
  ,("004237325","Ricardo",751.28)
  ,("004452476","Ivone",502.76)
  ,("004221454","Manuela",455.58)
  ,("005637820","Guilherme",439.12)
  ,("004508526","Cassio",400.9)
  ,("004211533","Valeria",375)
  ,("004458563","Luiz",278.42)
  ,("004207641","Magali",250)
  ,("004556853","Marcel",241.93)
  ,("008177213","Pedro",239.36)
  ,("001719494","Luis",235.55)
  ,("004261201","Ana",220)
  ,("004322719","Gisela",220)
  ,("003512801","Lais",211.2)
  ,("004359408","Hepta",187.01)
  ,("004377713","Danieli",164.55)
  ,("004328934","Valeria",163.84)
  ,("004405476","Mariana",98.98)
  ,("004907688","Heitor",97.41)
  ,("004431591","Mario",94.24)
  ,("004974089","Celia",93.23)
  ,("005701765","F",92.31)
  ,("004340223","Pedro",88.45)
  ,("004360431","Carlos",87.52)
  ,("004809902","Pedro",86.78)
  ,("004212132","Joao",86.38)
  ,("004384258","Paula",86.25)
  ,("001368670","Thiago",85.95)
  ,("000827730","Luciana",84.62)
  ,("005076418","Leonardo",81.87)
  ,("004230529","Lais",72.24)
  ,("004404342","Adson",54.9)
  ,("004400640","Felipe",51.44)
  ,("005142661","Sabrina",50.86)
  ,("004454365","Rafael",50.32)
  ,("001731007","Guilherme",44.79)
  ,("005245032","Rosa",39.91)
  ,("004238164","Daniela",39.15)
  ,("008336332","Carlos",37.15)
  ,("004806286","Vera",35.8)
  ,("004264780","Marcelo",35.29)
  ,("004119016","Hemat",35.24)
  ,("004340036","Eduardo",33.82)
  ,("004452597","Lara",31.21)
  ,("002894447","Joao",31.04)
  ,("004404724","Leandro",30.02)
  ,("005927101","Simone",30)
  ,("004377415","Angela",28.73)
  ,("004389994","Polyanna",26)
  ,("004350197","Gisela",25.08)
  ,("005186167","Andrea",22.15)
  ,("004371857","Nazareth",21.52)
  ,("004388077","Wladmir",20.89)
  ,("004214604","Maria",20.75)
  ,("004204255","Amado",18.77)
  ,("005374916","Marco",17.86)
  ,("004422594","Wandir",14.67)
  ,("005905713","Neila",13.77)
  ,("005135105","Brenner",11.26)
  ,("004976625","Norton",10.76)
  ,("004480134","Jose",10.61)
  ,("004419141","Paulo",10.42)
  ,("004420763","Christian",10.3)
  ,("004216298","Flordeliz",9.82)
  ,("001294033","Viviane",8.92)
  ,("008004995","Jose",8.2)
  ,("005043894","Naiara",8)
  ,("004530494","Rosangela",7.05)
  ,("005268516","Luis",6.35)
  ,("008013889","Carolina",6.24)
  ,("005905737","Cairo",6.17)
  ,("004448501","Joao",5.55)
  ,("004756968","Daniely",4.66)
  ,("008244362","Lincoln",3.28)
  ,("004308815","Zeli",2.51)
  ,("004419765","Walter",2.41)
  ,("004340217","Augusto",1.87)
  ,("001882235","Lago",1.84)
  ,("004840589","Leda",1.8)
  ,("004460491","Pedro",1.72)
  ,("005022526","Alexandre",1.7)
  ,("004220849","Dulce",1.68)
  ,("004335251","Edmundo",1.54)
  ,("005575050","Garcia",1.54)
  ,("004855960","Cleria",1.51)
  ,("004214460","Maria",1.5)
  ,("004713953","Alessandra",1.47)
  ,("004805333","Rosana",1.47)
  ,("004218542","Jose",1.45)
  ,("004527606","Marcia",1.43)
  ,("005886225","Vinicius",1.43)
  ,("004321092","Daniel",1.37)
  ,("005142624","Rodrigo",1.36)
  ,("004451652","Mateus",1.35)
  ,("004587511","Carlos",1.34)
  ,("004213139","Leonardo",1.3)
  ,("004381415","Joao",1.26)
  ,("004482090","Cezar",1.25)
  ,("005018038","Elaine",1.25)
  ,("004975924","Sergio",1.24)
  ,("008328804","Sonia",1.19)
  ,("004313254","Gustavo",1.17)
  ,("004911541","Tiago",1.16)
  ,("004267976","E3",1.12)
  ,("004479734","Rodrigo",1.12)
  ,("005381719","Maria",1.11)
  ,("004693308","Laura",1.03)
  ,("004462543","Rodolfo",1.01)
  ,("004360430","Viomar",1)
  ,("005440756","Valeria",1)
  ,("004486497","Elena",0.96)
  ,("004242237","Mariah",0.95)
  ,("005146441","Jose",0.95)
  ,("008012870","Ana",0.92)
  ,("005535788","Emilia",0.89)
  ,("004115403","Hebert",0.88)
  ,("001759765","Natal",0.86)
  ,("005428871","Rosangela",0.84)
  ,("004497875","Henrique",0.83)
  ,("005232019","Pedro",0.83)
  ,("005683532","Sylverson",0.82)
  ,("004223502","Bruna",0.78)
  ,("001000288","Isabella",0.73)
  ,("008298906","Fernando",0.72)
  ,("004346716","Tiago",0.71)
  ,("004425261","Thaysa",0.71)
  ,("004588677","Rachel",0.71)
  ,("005660155","Carolina",0.7)
  ,("004335144","Edmundo",0.67)
  ,("008004835","Sergio",0.67)
  ,("004473942","Daianne",0.62)
  ,("005924958","Tiago",0.6)
  ,("005338054","Elaine",0.57)
  ,("008119302","Vitor",0.56)
  ,("002786022","Paulo",0.55)
  ,("005105970","Vera",0.53)
  ,("005133039","Paulo",0.5)
  ,("004452946","Otavio",0.49)
  ,("005110894","Ana",0.49)
  ,("008072033","Leticia",0.49)
  ,("005558076","Alexandre",0.47)
  ,("008035153","Claudio",0.47)
  ,("003641655","Marcelo",0.45)
  ,("004805133","Patricia",0.45)
  ,("004767746","Isabele",0.44)
  ,("005662526","Aguinaldo",0.44)
  ,("005949170","Cintia",0.44)
  ,("005141215","Karina",0.42)
  ,("004432935","Jose",0.41)
  ,("005216881","Renan",0.41)
  ,("005341184","Breno",0.41)
  ,("005530256","Carolina",0.41)
  ,("004424671","Luisa",0.4)
  ,("004424761","Pedro",0.4)
  ,("004563252","Fernando",0.4)
  ,("005654767","Diego",0.4)
  ,("004453302","Isabella",0.39)
  ,("004466350","Raquel",0.39)
  ,("004397124","Murylo",0.38)
  ,("004424714","Helena",0.38)
  ,("002973105","Darlan",0.37)
  ,("004972070","Maria",0.37)
  ,("005725431","Bruno",0.37)
  ,("003836362","Isabella",0.33)
  ,("004556150","Marina",0.33)
  ,("004890544","Assako",0.3)
  ,("005379541","Jose",0.28)
  ,("004908680","Elene",0.26)
  ,("005270025","Denize",0.25)
  ,("008008723","Redrau",0.25)
  ,("004207278","Cesar",0.23)
  ,("005035754","Jose",0.23)
  ,("005092207","Bruno",0.23)
  ,("005110867","Dig",0.23)
  ,("004278033","Daisy",0.21)
  ,("004612043","Yuri",0.21)
  ,("004848843","Alaercio",0.21)
  ,("005697554","Maria",0.2)
  ,("005304669","Emilson",0.18)
  ,("004432455","Luciana",0.17)
  ,("002687737","Jose",0.16)
  ,("002694089","Vitor",0.16)
  ,("004339183","Jalison",0.15)
  ,("004357159","Joao",0.15)
  ,("004398253","Euler",0.15)
  ,("004754056","Bruno",0.15)
  ,("004805269","Clisia",0.15)
  ,("004320840","Natalia",0.14)
  ,("005075382","Nayara",0.14)
  ,("008336728","Arthur",0.14)
  ,("000938440","Base",0.12)
  ,("004382374","Theomar",0.12)
  ,("004551472","Diego",0.12)
  ,("005313179","Maira",0.12)
  ,("008071998","Isadora",0.11)
  ,("004281300","Franklin",0.1)
  ,("005880251","Luiz",0.1)
  ,("008351535","David",0.1)
  ,("004451996","Adriano",0.09)
  ,("004493324","Daniel",0.09)
  ,("005047946","Gabriel",0.09)
  ,("008407512","Mauricio",0.09)
  ,("004213373","Alexandre",0.08)
  ,("004216434","Jaime",0.08)
  ,("005720146","Jordanna",0.07)
  ,("004691225","Anna",0.06)
  ,("008026930","Joao",0.06)
  ,("008037529","Melissa",0.05)
  ,("008123677","Priscilla",0.05)
  ,("005274028","Rafael",0.04)
  ,("003107135","Ana",0.03)
  ,("004329229","Gabriel",0.03)
  ,("004946997","Eduardo",0.03)
  ,("008070544","Marina",0.03)
  ,("008110684","Edval",0.03)
  ,("008365179","Giovana",0.03)
  ,("004384131","Andre",0.02)
  ,("004457389","Rafael",0.02)
  ,("005295509","Bhruna",0.02)
  ,("008090243","Gabriel",0.02)
  ,("002878817","Guilherme",0.01)
  ,("004223226","Yeshua",0.01)
  ,("004272426","Rodrigo",0.01)
  ,("004400000","Vilma",0.01)
  ,("004462930","Walter",0.01)
  ,("004472076","Rubens",0.01)
  ,("004862746","Cesar",0.01)
  ,("005068961","Jorge",0.01)
  ,("005142592","Alberto",0.01)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = 2 + $i
  $row = $data[$i]
  $ws.Cells.Item($r, 1).NumberFormat = "@"
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = [double]$row[2]
}
